# Update the public EPEX spot prices workbook with the latest day of data.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Prix Spot": add a new date column CR ("17-sep") with 24 hourly
# price values, mirroring the format of the preceding column CQ.
# ---------------------------------------------------------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")

# Header cell: copy formatting (bold, border, centered) from CQ1.
$headerSrc = $wsSpot.Range("CQ1")
$headerDst = $wsSpot.Range("CR1")
$headerDst.Value = "17-sep"
$headerSrc.Copy()
$headerDst.PasteSpecial(-4122)
$excel.CutCopyMode = $false

$hourlyValues = @{
    2  = 71.56999999999999
    3  = 63.37
    4  = 50.98
    5  = 40.32
    6  = 32
    7  = 35.16
    8  = 67.78
    9  = 111.07
    10 = 112.08
    11 = 85.84999999999999
    12 = 57.9
    13 = 37.6
    14 = 18.05
    15 = 7.88
    16 = 0.09
    17 = 6.11
    18 = 12.63
    19 = 50
    20 = 83.12
    21 = 111.76
    22 = 123.55
    23 = 95.2
    24 = 87.36
    25 = 63.02
}

foreach ($row in 2..25) {
    $wsSpot.Cells.Item($row, 96).Value = $hourlyValues[$row]
}

# ---------------------------------------------------------------------
# Sheet "Gaz": append the latest daily gas price as row 93.
# The date is written as literal text (not an auto-converted date
# serial) by evaluating it as a formula and then pasting the result
# back in as a plain value, which keeps the default (unstyled) cell
# formatting intact.
# ---------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")
$gazDate = $wsGaz.Range("A93")
$gazDate.Formula = '="2025-09-15"'
$gazDate.Copy()
$gazDate.PasteSpecial(-4163)
$excel.CutCopyMode = $false
$wsGaz.Range("B93").Value = 31.55

# ---------------------------------------------------------------------
# Sheet "CO2": append the latest daily CO2 price as row 93.
# ---------------------------------------------------------------------
$wsCO2 = $wb.Worksheets.Item("CO2")
$co2Date = $wsCO2.Range("A93")
$co2Date.Formula = '="2025-09-15"'
$co2Date.Copy()
$co2Date.PasteSpecial(-4163)
$excel.CutCopyMode = $false
$wsCO2.Range("B93").Value = 76.23999999999999
